$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drawing: remove the per-equipment placeholder pictures (Image 3..6), ---
# --- keep the two corporate logo pictures (image1.png / image2.png).      ---
$ws.Shapes.Item("Image 3").Delete()
$ws.Shapes.Item("Image 4").Delete()
$ws.Shapes.Item("Image 5").Delete()
$ws.Shapes.Item("Image 6").Delete()

# --- Row 18: clear out the sample history entry but keep the row/merges ---
$ws.Range("A18").ClearContents()
$ws.Range("B18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("G18").ClearContents()
$ws.Range("H18").ClearContents()

# Helper style source: the blank template cells on row 16 (plain, borderless).
$ws.Range("B16").Copy()

# --- Row 19: wipe the sample data back into a blank template row ---
$ws.Range("B19:D19").UnMerge()
$ws.Range("E19:F19").UnMerge()
$ws.Range("J19:K19").UnMerge()
$ws.Range("A19:K19").Clear()
"A19","B19","C19","E19","F19","G19","H19","I19","J19","K19" | ForEach-Object {
    $ws.Range($_).PasteSpecial(-4122)
}
$ws.Rows.Item(19).RowHeight = 12.75

# --- Row 20: same treatment as row 19 ---
$ws.Range("B20:D20").UnMerge()
$ws.Range("E20:F20").UnMerge()
$ws.Range("J20:K20").UnMerge()
$ws.Range("A20:K20").Clear()
"A20","B20","C20","E20","F20","G20","H20","I20","J20","K20" | ForEach-Object {
    $ws.Range($_).PasteSpecial(-4122)
}
$ws.Rows.Item(20).RowHeight = 12.75

# --- Row 21: same treatment, but this one keeps its D cell ---
$ws.Range("B21:D21").UnMerge()
$ws.Range("E21:F21").UnMerge()
$ws.Range("J21:K21").UnMerge()
$ws.Range("A21:K21").Clear()
"A21","B21","C21","D21","E21","F21","G21","H21","I21","J21","K21" | ForEach-Object {
    $ws.Range($_).PasteSpecial(-4122)
}
$ws.Rows.Item(21).RowHeight = 12.75

# --- Row 22: fully blank (no cells at all) ---
$ws.Range("B22:D22").UnMerge()
$ws.Range("E22:F22").UnMerge()
$ws.Range("J22:K22").UnMerge()
$ws.Range("A22:K22").Clear()
$ws.Rows.Item(22).RowHeight = 15.75
